$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16
$rng = $ws.Range("A$row" + ":H$row")
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

$ws.Range("A$row").Value = "2025-08-15 09:40:31 UTC"
$ws.Range("B$row").Value = "2025-08-15 15:10:31 IST"
$ws.Range("C$row").Value = "SKIPPED"
$ws.Range("D$row").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E$row").Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("F$row").Value = ""
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = ""
